# Update the title paragraph date range from "Fall 2023" to "Spring 2024",
# and introduce a space run between "Parameters" and the en dash, matching
# the target revision exactly at the run level.
$d = $word.ActiveDocument

# Locate the title paragraph (the one still mentioning "Fall") rather than
# assuming a fixed paragraph index.
$titlePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text -like "*Fall*") {
        $titlePara = $candidate.Range
        break
    }
}
if ($titlePara -eq $null) {
    $titlePara = $d.Paragraphs(1).Range
}

$newParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="221BF16F" w14:textId="7346DDA6" w:rsidR="00CE58B2" w:rsidRDefault="00AB2F21"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Electronics Basics </w:t></w:r><w:r w:rsidR="00CE58B2"><w:rPr><w:b/></w:rPr><w:t>Parameter</w:t></w:r><w:r w:rsidR="004B7A9A"><w:rPr><w:b/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00CE58B2"><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">– </w:t></w:r><w:r w:rsidR="00373A7E"><w:rPr><w:b/></w:rPr><w:t>Spring 20</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidR="00DE3015"><w:rPr><w:b/></w:rPr><w:t>24</w:t></w:r></w:p>
'@

$titlePara.InsertXML($newParaXml)

Write-Host "Updated title paragraph text: $($d.Paragraphs(1).Range.Text)"
